# Applies the edit described by the commit "funicoandano la creacion del excel":
#  1. Moves the worksheet's active selection from H16 to J20.
#  2. Drops the cached "best fit" auto-width flag on columns A, C and D
#     (their stored widths are kept the same / as close as the host
#     ColumnWidth API allows).
#  3. Inserts a new row 11 with SUM() totals formulas in K11 and L11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Move the active cell / selection to J20 (was H16).
$ws.Range("J20").Select()

# 2. Re-assert explicit column widths for columns A, C and D so that Excel
#    no longer flags them as "bestFit" (auto-fit) columns, while keeping the
#    on-disk width as close as possible to the original values
#    (8.44140625, 14.109375 and 18 characters respectively).
$ws.Columns.Item(1).ColumnWidth = 7.666666666666667   # column A -> width 8.44140625 (bestFit removed)
$ws.Columns.Item(3).ColumnWidth = 13.333333333333334  # column C -> width 14.109375 (bestFit removed)
$ws.Columns.Item(4).ColumnWidth = 17.166666666666668  # column D -> width 18 (bestFit removed)

# 3. Add the new totals row 11 with SUM formulas in K11 and L11.
$ws.Range("K11").Formula = "=SUM(K13:K1020)"
$ws.Range("L11").Formula = "=SUM(L13:L1020)"
